$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows at row 81 (shifts rows 81-93 down to 83-95).
$ws.Rows("81:82").Insert(-4121, 0)

# The newly inserted rows inherit odd auto-generated styles; fix them by
# copying the format from row 80 (which has the exact style pattern we need).
$ws.Range("A80:G80").Copy()
$ws.Range("A81:G82").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new row 81 with the new resource-bundle entry.
$ws.Range("B81").Value = "XML2SOURCE_FILE.ROUTECONFIG.LIST"
$ws.Range("C81").Value = "ページコンポーネントをvue-routerでロードするための設定ファイルです"

# Row 82 stays blank (already cleared by the insert).
